$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray empty-string cells on row 74 (I74, K74:R74) so they become
# truly empty (not present) as in the target.
$ws.Range("I74").ClearContents()
$ws.Range("K74:R74").ClearContents()

# New data rows 75-89. Last field notes whether column I should be written
# as an explicit empty (quote-prefixed) text cell, matching rows 80-89 in
# the target sheet.
$data = @(
    @("6XS17353","TOPICREM CREMA CALMANTE 40ML","ANEXOS","No Tiene PT - TRADUZIDO","Tiene ES","Tiene IT","40","ML","Revisado y Traducido",$false),
    @("6XS18000","TOPICREM BABY 2EN1 GEL LIMPIADOR 500ML","ANEXOS","No Tiene PT - TRADUZIDO","Tiene ES","Tiene IT","500","ML","Revisado y Traducido",$false),
    @("2CC05367","KRISTIN THE ONE PURPLE SHAMPOO 296 ML","CABELLO ACONDIC. SUAVIZANTE","No Tiene PT - TRADUZIDO","Tiene ES","Tiene IT","296","ML","Revisado y Traducido",$false),
    @("4EF05522","BENCH TOGETHER SET HER 30 ML+75ML SG","PERF. ESTUCHES MUJER","No Tiene PT - TRADUZIDO","Tiene ES","No Tiene IT - TRADOTTO","2","UND","Revisado y Traducido",$false),
    @("4EM03368","BENCH TOGETHER SET HIM 30 ML+75ML SG","PERF. ESTUCHES HOMBRE","No Tiene PT - TRADUZIDO","Tiene ES","No Tiene IT - TRADOTTO","2","UND","Revisado y Traducido",$false),
    @("2CC05963","INOPHARM RENU RECONSTRUCTING CHAMPU 250ML","VARIOS","No Tiene PT - TRADUZIDO","Tiene ES","No Tiene IT - TRADOTTO","1","UND","Revisado y Traducido",$true),
    @("2CC05960","INOPHARM SCALP CARE ANTI-DANDRUFF CHAMPU 250ML","VARIOS","No Tiene PT - TRADUZIDO","Tiene ES","No Tiene IT - TRADOTTO","1","UND","Revisado y Traducido",$true),
    @("2CC05961","INOPHARM SCALP CARE MICELLAR ANTI-DANDRUFF 250ML","VARIOS","No Tiene PT - TRADUZIDO","Tiene ES","No Tiene IT - TRADOTTO","1","UND","Revisado y Traducido",$true),
    @("2CC05962","INOPHARM SCALP CARE REGROWTH CHAMPU 250ML","VARIOS","No Tiene PT - TRADUZIDO","Tiene ES","No Tiene IT - TRADOTTO","1","UND","Revisado y Traducido",$true),
    @("6XS18637","INOPHARM INTIMA EVERYDAY USE INTIMATE GEL 250ML","VARIOS","No Tiene PT - TRADUZIDO","Tiene ES","No Tiene IT - TRADOTTO","1","UND","Revisado y Traducido",$true),
    @("6XS18639","INOPHARM INTIMA EXTRA PROTECTION CARE 250ML","VARIOS","No Tiene PT - TRADUZIDO","Tiene ES","No Tiene IT - TRADOTTO","1","UND","Revisado y Traducido",$true),
    @("6XS18640","INOPHARM INTIMA INTIMATE GEL FEMENINE GUARD 250ML","VARIOS","No Tiene PT - TRADUZIDO","Tiene ES","No Tiene IT - TRADOTTO","1","UND","Revisado y Traducido",$true),
    @("6XS18641","INOPHARM INTIMA INTIMATE OIL MENOPAUSE 250ML","VARIOS","No Tiene PT - TRADUZIDO","Tiene ES","No Tiene IT - TRADOTTO","1","UND","Revisado y Traducido",$true),
    @("2CA06703","INOPHARM RENU RECONSTRUCTING CONDITONER 250ML","VARIOS","No Tiene PT - TRADUZIDO","Tiene ES","No Tiene IT - TRADOTTO","1","UND","Revisado y Traducido",$true),
    @("2CA06704","INOPHARM RENU RECONSTRUCTING HAIR MASK 250ML","VARIOS","No Tiene PT - TRADUZIDO","Tiene ES","No Tiene IT - TRADOTTO","1","UND","Revisado y Traducido",$true)
)

$row = 75
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
    # Quantity column: source data stores this as text (matches the rest of
    # the column throughout the sheet), so force text with a leading
    # apostrophe rather than letting it be auto-detected as a number.
    $ws.Cells.Item($row, 7).Value = "'" + $item[6]
    $ws.Cells.Item($row, 8).Value = $item[7]
    if ($item[9]) {
        # Explicit empty (quote-prefixed) text cell, matching the target's
        # <c t="inlineStr"/> placeholder cells on rows 80-89.
        $ws.Cells.Item($row, 9).Value = "'"
    }
    $ws.Cells.Item($row, 10).Value = $item[8]
    $row++
}
